$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format first so numeric-looking strings
# (e.g. "1.00", "0.198", "0.0000184") are not silently reinterpreted by Excel
# as numbers, which would change their literal display text.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '92.902.67'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '3.428.13'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '232.12'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').Value = '620.85'
$ws.Range('E6').Value = '  -3.44%  '
$ws.Range('D7').Value = '1.38'
$ws.Range('E7').Value = '  -4.86%  '
$ws.Range('D8').Value = '0.392'
$ws.Range('E8').Value = '  -3.79%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').Value = '3.427.47'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '42.96'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').Value = '0.198'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '6.26'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '4.065.48'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '92.809.77'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').Value = '0.0000246'
$ws.Range('E17').Value = '  -2.22%  '
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').Value = '3.442.10'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = '17.85'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '11.63'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '498.44'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = '3.32'
$ws.Range('E23').Value = '  +2.54%  '
$ws.Range('D24').Value = '0.441'
$ws.Range('E24').Value = '  -12.21%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000184'
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '6.34'
$ws.Range('E26').Value = '  -2.76%  '
$ws.Range('D27').Value = '91.21'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('D28').Value = '12.00'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = '3.598.22'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '11.32'
$ws.Range('E30').Value = '  -3.28%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value = '1.01'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.72'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.135'
$ws.Range('E33').Value = '  -2.91%  '
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').Value = '0.173'
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('D36').Value = '29.66'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').Value = '0.541'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').Value = '558.88'
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('D39').Value = '7.50'
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '1.39'
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.150'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').Value = '0.920'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '1.73'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '23.68'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').Value = '3.68'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '5.48'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0407'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '53.04'
$ws.Range('E49').Value = '  -3.72%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '2.11'
$ws.Range('E50').Value = '  -3.96%  '
$ws.Range('B51').Value = 'Fantom'
$ws.Range('C51').Value = 'https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm'
$ws.Range('D51').Value = '1.12'
$ws.Range('E51').Value = '  +17.34%  '

# Remove the temporary text-format styling so the cells end up with the same
# (default) style they started with, leaving only the text content changed.
$dataRange.ClearFormats()
